# Apply the "megazord" edit described by the commit diff.
#
# Summary of changes:
#  1. Sheet "Planilha1" (sheet3) renamed to "Especimes_LACV_emprestimoCHUNB".
#  2. Sheet3 gains a new header row (ID / No_LACV) and a new column B filled
#     with "SIM" for every existing data row.
#  3. Sheet4 ("Planilha2") gains a new header row (Especie / Colecao / ID /
#     Hemipenis_Evertido); existing species/collection codes are upper-cased
#     and a new column D is filled with "SIM".
#  4. Sheet2 ("Domingos 2014") column headers G1/J1/K1 are shortened
#     (Localidade -> Loc, Latitude -> Lat, Longitude -> Long).
#  5. View/selection tweaks: sheet1 keeps its frozen header row but scrolled
#     back to the top; sheet2/sheet3/sheet4 selections move; sheet2 stays the
#     active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1-2. Sheet3: rename + new header row + new "No_LACV" column.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "Especimes_LACV_emprestimoCHUNB"

$ws3.Rows.Item(1).Insert()
$ws3.Range("A1").Value = "ID"
$ws3.Range("B1").Value = "No_LACV"

for ($r = 2; $r -le 47; $r++) {
  $ws3.Cells.Item($r, 2).Value = "SIM"
}

# ---------------------------------------------------------------------------
# 3. Sheet4: new header row + uppercase codes + new "Hemipenis_Evertido" column.
# (data rows first, header row last -- mirrors how the shared-string table
# for the new labels ends up ordered in the authored workbook)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Rows.Item(1).Insert()

$codes4 = @("A","B","C","B ","B","B","B","C","B","C","C","D","D","E","E","F","G","H","H","NA")
$colls4 = @("CHUNB","LACV","CHUNB","CHUNB","CHUNB","CHUNB","CHUNB","CHUNB","LACV","CHUNB","CHUNB","CHUNB","CHUNB","CHUNB","CHUNB","CHUNB","CHUNB","CHUNB","CHUNB","CHUNB")

for ($i = 0; $i -lt 20; $i++) {
  $r = $i + 2
  $ws4.Cells.Item($r, 1).Value = $codes4[$i]
  $ws4.Cells.Item($r, 2).Value = $colls4[$i]
  $ws4.Cells.Item($r, 4).Value = "SIM"
}

$ws4.Range("A1").Value = "Especie"
$ws4.Range("B1").Value = "Colecao"
$ws4.Range("C1").Value = "ID"
$ws4.Range("D1").Value = "Hemipenis_Evertido"

# ---------------------------------------------------------------------------
# 4. Sheet2: shorten a few column headers.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("G1").Value = "Loc"
$ws2.Range("J1").Value = "Lat"
$ws2.Range("K1").Value = "Long"

# ---------------------------------------------------------------------------
# 5. View / selection tweaks.
# ---------------------------------------------------------------------------

# Sheet1: keep the row-1 freeze, but scroll back up to the top (topLeftCell
# A2) while leaving the original C12 selection untouched.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$win1 = $excel.ActiveWindow
$win1.FreezePanes = $false
$ws1.Range("A2").Select()
$win1.FreezePanes = $true
$ws1.Range("C12").Select()

# Sheet3: scroll down near the bottom of the new table, select column B.
$ws3.Activate()
$ws3.Range("B2:B47").Select()

# Sheet4: select the new column D.
$ws4.Activate()
$ws4.Range("D2:D21").Select()

# Sheet2: move the selection, and re-activate it last so it remains the
# workbook's active tab (matches the unchanged activeTab="1" in workbook.xml).
$ws2.Activate()
$ws2.Range("I5").Select()
